$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "MEDINA VALLEJOS ERICK LEONARDO",
    "ZAMORA TAMAY NEYSER IVAN",
    "SOTO VILLENA NILSON",
    "VASQUEZ DIAZ LUZ ANGELICA",
    "ROJAS VASQUEZ FLOR NOELITA",
    "PÓSITO CHUGDEN NANIX",
    "SOTO VALLEJOS ELSITA",
    "TIRADO PEREZ JEINER",
    "TELLO FERNANDEZ MILENY",
    "VASQUEZ LUNA YUDITH",
    "GALLARDO CORTEZ MELISSA DEL CARMEN",
    "RUIZ RUIZ LUZ MERI",
    "BENAVIDES MARRUFO ARACELI",
    "BENAVIDES SALAZAR IDELSA"
)

$totals = @(84, 83, 81, 80, 79, 79, 78, 77, 75, 73, 71, 69, 66, 66)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}
